# Add 2022-Q1 data
# 1) Create the new "2022-Q1" worksheet (same layout as the other quarter
#    sheets) positioned right before the "总计" (total) sheet.
# 2) Fill it with the 2022-Q1 fund holdings.
# 3) Insert a new summary row into "总计" for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create "2022-Q1" by copying the "2021-Q4" sheet (identical
# header/format) and dropping it in right after "2021-Q4" (i.e. right
# before "总计").
# ---------------------------------------------------------------------
$srcQuarter = $wb.Worksheets.Item("2021-Q4")
$srcQuarter.Copy($null, $srcQuarter)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Remove the rows copied from 2021-Q4 that we don't need (rows 13-36),
# keeping the header (row 1) and 11 data rows (rows 2-12).
$q1.Range("A13:A36").EntireRow.Delete()

# ---------------------------------------------------------------------
# Step 2: populate the 2022-Q1 holdings data.
# Columns: A index, B code, C name, D size, E stock-position-total,
#          F position-share, G holding value (yi), H position rank
# ---------------------------------------------------------------------
$data = @(
    @("000362", "国泰聚信价值优势灵活配置混合A", "56.15", "89.00", "3.52", "1.9765", 10),
    @("012173", "国泰兴泽优选一年持有期混合A",   "16.89", "89.83", "3.57", "0.6030", 8),
    @("000363", "国泰聚信价值优势灵活配置混合C", "17.09", "89.00", "3.52", "0.6016", 10),
    @("011230", "创金合信数字经济主题股票C",     "17.18", "92.17", "3.24", "0.5566", 8),
    @("011229", "创金合信数字经济主题股票A",     "12.18", "92.17", "3.24", "0.3946", 8),
    @("012174", "国泰兴泽优选一年持有期混合C",   "7.14",  "89.83", "3.57", "0.2549", 8),
    @("003713", "英大睿盛灵活配置混合A",         "5.99",  "87.42", "3.69", "0.2210", 10),
    @("003714", "英大睿盛灵活配置混合C",         "2.40",  "87.42", "3.69", "0.0886", 10),
    @("014339", "长江智能制造混合A",             "3.28",  "21.63", "0.80", "0.0262", 10),
    @("001608", "英大策略优选混合C",             "0.03",  "89.86", "4.48", "0.0013", 9),
    @("014340", "长江智能制造混合C",             "0.15",  "21.63", "0.80", "0.0012", 10)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $q1.Range("A$row").Value = $i

    # Columns that must stay text (fund codes / decimal-looking text
    # fields) need an explicit text format so Excel doesn't coerce them
    # into numbers (which would drop leading zeros / trailing zeros).
    $textRange = $q1.Range("B$row`:G$row")
    $textRange.NumberFormat = "@"

    $q1.Range("B$row").Value = $rec[0]
    $q1.Range("C$row").Value = $rec[1]
    $q1.Range("D$row").Value = $rec[2]
    $q1.Range("E$row").Value = $rec[3]
    $q1.Range("F$row").Value = $rec[4]
    $q1.Range("G$row").Value = $rec[5]

    # Drop the NumberFormat residue left behind on B:G so the cells keep
    # the default (no explicit style index), matching the sibling rows.
    $textRange.Style = "Normal"

    $q1.Range("H$row").Value = $rec[6]
}

# ---------------------------------------------------------------------
# Step 3: add the 2022-Q1 row to the "总计" (totals) summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Copy the format of the existing first data row (row 3, after the
# insert shifts it down) onto the newly inserted row 2, then fill in the
# values - this keeps the exact style indices used by the sibling rows.
$total.Range("A3:D3").Copy()
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 11
$total.Range("D2").Value = 4.73
